$d = $word.ActiveDocument

# The cover page originally reads "NIM 1215002" built from three runs:
#   "NIM" + " 121500" + "2"
# The edit replaces the middle run's leading space with a period, so the
# line reads "NIM.1215002" ("NIM" + "." + "121500" + "2").
$d.Content.Find.Execute(" 121500", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ".121500", 2)
